# Alteração nos rótulos da tabela para já transformar a primeira linha em
# cabeçalho automaticamente no Power BI.
# Prefix the year/interval header labels in row 1 of each sheet with
# "Ano " (for single years) or "Intervalo " (for year ranges / the
# "Potencia Incremental" sheet).

$wb = $excel.ActiveWorkbook

# Sheets whose header row (B1:E1) holds plain years -> prefix with "Ano "
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Text
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Text
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Text
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Text
}

# "Potencia Incremental" sheet uses year intervals -> prefix with "Intervalo "
$wsInc = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsInc.Range("B1").Value = "Intervalo " + $wsInc.Range("B1").Text
$wsInc.Range("C1").Value = "Intervalo " + $wsInc.Range("C1").Text
$wsInc.Range("D1").Value = "Intervalo " + $wsInc.Range("D1").Text
$wsInc.Range("E1").Value = "Intervalo " + $wsInc.Range("E1").Text

# "Custo Total" sheet only has a single year column (B1) -> prefix with "Ano "
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano " + $wsCusto.Range("B1").Text
